# HMB_ScintRdoutReport.docx edit
#
# 1) Merge the three runs "can transmit and rece" / "ive 32-bit words
#    simultaneously " / "with another " (which were split apart by a
#    "_GoBack" bookmark sitting between them) back into a single run,
#    dropping the stray bookmark in the process.
# 2) The document's lone comment ("Fact check" by Kevin) is renumbered
#    from w:id="1" down to w:id="0" (as Word does when comment ids get
#    compacted). Recreate the comment anchored on the same range with
#    the same author/initials so the commentRangeStart/End/Reference and
#    comments.xml all come out with id 0.
# 3) Split the trailing " were successful  " run so the two trailing
#    spaces are pushed into their own run, with a (moved) "_GoBack"
#    bookmark marking the former cursor position right after
#    "successful".

$d = $word.ActiveDocument

# --- Change 1: stitch "can transmit and rece" + "ive 32-bit words
#     simultaneously " + "with another " into one run, removing the
#     "_GoBack" bookmark that used to sit between them. Find/Replace
#     across the run boundary naturally merges the text into the
#     formatting of the run under the cursor and swallows the bookmark. ---
$d.Content.Find.Execute(
    "can transmit and receive 32-bit words simultaneously with another",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can transmit and receive 32-bit words simultaneously with another", 2) | Out-Null

# --- Change 2: renumber the comment from id 1 to id 0. The object model
#     doesn't expose comment ids directly, but a freshly mustered comment
#     collection numbers its first member 0, so delete the existing
#     comment and re-add an equivalent one anchored on the same range. ---
$comments = $d.Comments
$oldComment = $comments.Item(1)
$scope = $oldComment.Scope
$rangeStart = $scope.Start
$rangeEnd = $scope.End
$commentText = $oldComment.Range.Text

$oldComment.Delete()

$anchorRange = $d.Range($rangeStart, $rangeEnd)
$comments.Add($anchorRange, $commentText) | Out-Null

$newComment = $comments.Item(1)
$newComment.Author = "Kevin"
$newComment.Initial = "K"

# --- Change 3: split off the two trailing spaces after "successful"
#     into their own run and drop a (relocated) "_GoBack" bookmark right
#     before them, mirroring Word leaving its "last edit" bookmark at the
#     spot that was just typed. ---
$tail = $d.Content
$tail.Find.Execute("successful", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Collapse(0)
$d.Bookmarks.Add("_GoBack", $tail) | Out-Null
